# The deck currently carries two themes:
#   ppt/theme/theme1.xml -> "Office Theme" (used only by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"     (used by the Slide Master / the
#                                            actual slides, and by
#                                            presentation.xml's rId1)
#
# The target edit swaps the two themes' content: theme2.xml (the theme that
# drives every slide's look) must end up holding the plain "Office Theme"
# 12-colour scheme that theme1.xml currently has, while theme1.xml would take
# on the "Integral" colours. The font scheme (fontScheme) and the fill/line/
# effect scheme (fmtScheme) are byte-identical between the two themes, so the
# whole swap is really just the 12 theme colours (+ cosmetic "name"
# attributes, which PowerPoint's object model does not expose a setter for).
#
# PowerPoint's object model only exposes the slide-facing theme's colours
# (Slide.ThemeColorScheme / SlideRange.ThemeColorScheme) for editing, so we
# recolor that one (backed by ppt/theme/theme2.xml) to the "Office Theme"
# palette here.

$p = $ppt.ActivePresentation

function ConvertTo-OleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target palette = the stock "Office Theme" colour scheme (12 slots, in the
# fixed clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColor = $themeColors.Colors($i)
    $themeColor.RGB = ConvertTo-OleColor $officeThemeColors[$i - 1]
}
